$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.128
$ws.Range("E2").Value = 0.161
$ws.Range("F2").Value = -0.0211
$ws.Range("G2").Value = 0.2846318727879966
$ws.Range("H2").Value = 0.2787498130701361
$ws.Range("I2").Value = 0.2290340640891652
$ws.Range("J2").Value = 0.181932065945081
$ws.Range("K2").Value = 236.3
$ws.Range("L2").Value = 0.1177907382483426
$ws.Range("M2").Value = 203.9647
$ws.Range("N2").Value = 0.0520066039419669
$ws.Range("O2").Value = 0.8631599661447312
$ws.Range("P2").Value = 135.4427
$ws.Range("Q2").Value = 0.03453497029501007
$ws.Range("R2").Value = 0.5731811256876851
$ws.Range("S2").Value = 68.52199999999999
$ws.Range("T2").Value = 0.3359502894373389
$ws.Range("U2").Value = 1744.2
$ws.Range("V2").Value = 0.4447334200260078
$ws.Range("W2").Value = 0.1053369667400021
$ws.Range("X2").Value = 0.03182182098107456
$ws.Range("Y2").Value = 0.07351514575892751
$ws.Range("Z2").Value = 0.1636313741420596
$ws.Range("AA2").Value = 0.002132674291460412
$ws.Range("AB2").Value = 0.02795785722916758
$ws.Range("AC2").Value = -0.02531224957115663
$ws.Range("AD2").Value = 7857.3
$ws.Range("AE2").Value = 117.6738201536282
$ws.Range("AF2").Value = 7974.973820153628
$ws.Range("AG2").Value = 6230.773820153629
$ws.Range("AH2").Value = 0.6703419688829351
$ws.Range("AI2").Value = 0.6835043207967004
$ws.Range("AJ2").Value = 0.6137076725330417
$ws.Range("AK2").Value = 0.6278759984129556
$ws.Range("AL2").Value = 315.136
$ws.Range("AM2").Value = 288.576
$ws.Range("AN2").Value = 16.10762607626076
$ws.Range("AO2").Value = 1.432397441104793
$ws.Range("AP2").Value = 12.77321406345557
$ws.Range("AQ2").Value = 1.564232645819472

# Row 3
$ws.Range("D3").Value = 0.07290000000000001
$ws.Range("E3").Value = 0.226
$ws.Range("F3").Value = 0.0751
$ws.Range("G3").Value = 0.6188769414575865
$ws.Range("H3").Value = 0.5483870967741935
$ws.Range("I3").Value = 0.525089605734767
$ws.Range("J3").Value = 0.3628708587159464
$ws.Range("K3").Value = 67
$ws.Range("L3").Value = 0.4002389486260454
$ws.Range("M3").Value = 48.6
$ws.Range("N3").Value = 0.03458090223423936
$ws.Range("O3").Value = 0.7253731343283583
$ws.Range("P3").Value = 48.6
$ws.Range("Q3").Value = 0.03458090223423936
$ws.Range("R3").Value = 0.7253731343283583
$ws.Range("U3").Value = 159.2
$ws.Range("V3").Value = 0.1132773587590721
$ws.Range("W3").Value = 0.2019897497738921
$ws.Range("X3").Value = 0.01965885483020477
$ws.Range("Y3").Value = 0.1823308949436873
$ws.Range("Z3").Value = 3.938823529411765
$ws.Range("AA3").Value = 1.429284276448222
$ws.Range("AB3").Value = 0.01969341362710942
$ws.Range("AC3").Value = 1.409590862821112
$ws.Range("AD3").Value = 13.4
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 13.4
$ws.Range("AG3").Value = -145.8
$ws.Range("AH3").Value = 0.009444601071327881
$ws.Range("AI3").Value = 0.0391812865497076
$ws.Range("AJ3").Value = -0.1157510320736742
$ws.Range("AK3").Value = -0.7975929978118159
$ws.Range("AL3").Value = 0.536
$ws.Range("AM3").Value = -7.423999999999999
$ws.Range("AN3").Value = 0.1490545050055617
$ws.Range("AO3").Value = 163.9925373134328
$ws.Range("AP3").Value = -1.621802002224694
$ws.Range("AQ3").Value = -11.83997844827586

# Row 4
$ws.Range("B4").Value = "Unifin Financiera, S. A. B. de C. V. (BMV:UNIFIN A)"
$ws.Range("D4").Value = 0.135
$ws.Range("E4").Value = 0.0973
$ws.Range("G4").Value = 0.9184515621929652
$ws.Range("H4").Value = 0.9184515621929652
$ws.Range("I4").Value = 0.7142857142857143
$ws.Range("J4").Value = 0.5728429985855729
$ws.Range("K4").Value = 72.90000000000001
$ws.Range("L4").Value = 0.1432501473766948
$ws.Range("M4").Value = 125.2
$ws.Range("N4").Value = 0.1994901210962396
$ws.Range("O4").Value = 1.717421124828532
$ws.Range("P4").Value = 65.2
$ws.Range("Q4").Value = 0.1038878266411727
$ws.Range("R4").Value = 0.8943758573388203
$ws.Range("S4").Value = 60
$ws.Range("T4").Value = 0.4792332268370607
$ws.Range("U4").Value = 136
$ws.Range("V4").Value = 0.2166985340981517
$ws.Range("W4").Value = 0.1467391304347826
$ws.Range("X4").Value = 0.055631039971481
$ws.Range("Y4").Value = 0.09110809046330162
$ws.Range("Z4").Value = 0.1567582552981764
$ws.Range("AA4").Value = 0.08979786901805016
$ws.Range("AB4").Value = 0.0427819487316305
$ws.Range("AC4").Value = 0.04701592028641966
$ws.Range("AD4").Value = 3140.2
$ws.Range("AF4").Value = 3140.2
$ws.Range("AG4").Value = 3004.2
$ws.Range("AH4").Value = 0.8334306491852009
$ws.Range("AI4").Value = 0.8463466566045872
$ws.Range("AJ4").Value = 0.8271931273748555
$ws.Range("AK4").Value = 0.8405002378088018
$ws.Range("AL4").Value = 314.6
$ws.Range("AM4").Value = 296
$ws.Range("AN4").Value = 8.572754572754572
$ws.Range("AO4").Value = 1.155435473617292
$ws.Range("AP4").Value = 8.2014742014742
$ws.Range("AQ4").Value = 1.228040540540541

# Row 5
$ws.Range("D5").Value = 0.121
$ws.Range("E5").Value = 0.161
$ws.Range("K5").Value = 67.59999999999999
$ws.Range("L5").Value = 0.1695510408828693
$ws.Range("M5").Value = 4.76
$ws.Range("N5").Value = 0.0146551724137931
$ws.Range("O5").Value = 0.07041420118343196
$ws.Range("P5").Value = 3.17
$ws.Range("Q5").Value = 0.009759852216748768
$ws.Range("R5").Value = 0.04689349112426036
$ws.Range("S5").Value = 1.59
$ws.Range("T5").Value = 0.3340336134453781
$ws.Range("U5").Value = 541.2
$ws.Range("V5").Value = 1.666256157635468
$ws.Range("W5").Value = 0.1396694214876033
$ws.Range("X5").Value = 0.02271271159162502
$ws.Range("Y5").Value = 0.1169567098959783
$ws.Range("Z5").Value = 0.3572260550129917
$ws.Range("AB5").Value = 0.02517769419386753
$ws.Range("AC5").Value = -0.02517769419386753
$ws.Range("AD5").Value = 140.8
$ws.Range("AF5").Value = 140.8
$ws.Range("AG5").Value = -400.4
$ws.Range("AH5").Value = 0.302405498281787
$ws.Range("AI5").Value = 0.2147323471099588
$ws.Range("AJ5").Value = 5.296296296296295
$ws.Range("AK5").Value = -3.49694323144105

# Row 6
$ws.Range("D6").Value = -0.0286
$ws.Range("E6").ClearContents()
$ws.Range("F6").Value = -0.0211
$ws.Range("I6").Value = 0.01524331122523977
$ws.Range("J6").Value = 0.01524331122523977
$ws.Range("K6").Value = -49.2
$ws.Range("L6").Value = -0.09298809298809299
$ws.Range("M6").Value = 1.92
$ws.Range("N6").Value = 0.002473270642792735
$ws.Range("O6").Value = -0.03902439024390243
$ws.Range("P6").Value = -0
$ws.Range("Q6").Value = -0
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = 1.92
$ws.Range("T6").Value = 1
$ws.Range("U6").Value = 474.7
$ws.Range("V6").Value = 0.6114904031946413
$ws.Range("W6").Value = -0.04869840641393646
$ws.Range("X6").Value = 0.03302004995999393
$ws.Range("Y6").Value = -0.08171845637393038
$ws.Range("Z6").Value = 0.2798177193848991
$ws.Range("AA6").Value = 0.004265348582920823
$ws.Range("AB6").Value = 0.02971215353136656
$ws.Range("AC6").Value = -0.02544680494844574
$ws.Range("AD6").Value = 1329.7
$ws.Range("AE6").Value = 117.6738201536282
$ws.Range("AF6").Value = 1447.373820153628
$ws.Range("AG6").Value = 972.6738201536282
$ws.Range("AH6").Value = 0.6508930433212694
$ws.Range("AI6").Value = 0.5884652895123202
$ws.Range("AJ6").Value = 0.5561397254466561
$ws.Range("AK6").Value = 0.4900431504902129
$ws.Range("AN6").Value = 42.07911392405063
$ws.Range("AP6").Value = 30.78081709346925

# Row 7
$ws.Range("D7").Value = 0.154
$ws.Range("E7").Value = 0.176
$ws.Range("K7").Value = 29.9
$ws.Range("L7").Value = 0.1931524547803617
$ws.Range("M7").Value = 2.7947
$ws.Range("N7").Value = 0.005013814137064944
$ws.Range("O7").Value = 0.09346822742474915
$ws.Range("P7").Value = 2.7727
$ws.Range("Q7").Value = 0.004974345174022246
$ws.Range("R7").Value = 0.09273244147157191
$ws.Range("S7").Value = 0.0219999999999998
$ws.Range("T7").Value = 0.007872043510931334
$ws.Range("U7").Value = 415
$ws.Range("V7").Value = 0.7445281664872623
$ws.Range("W7").Value = 0.07100451199240085
$ws.Range("X7").Value = 0.0306235920021552
$ws.Range("Y7").Value = 0.04038091999024565
$ws.Range("Z7").Value = 0.0466574235939478
$ws.Range("AB7").Value = 0.0262035609269686
$ws.Range("AC7").Value = -0.0262035609269686
$ws.Range("AD7").Value = 853.8
$ws.Range("AF7").Value = 853.8
$ws.Range("AG7").Value = 438.8
$ws.Range("AH7").Value = 0.6050170068027211
$ws.Range("AI7").Value = 0.64903078677309
$ws.Range("AJ7").Value = 0.4404738004416784
$ws.Range("AK7").Value = 0.4872848417545808

# Row 8
$ws.Range("D8").Value = 0.143
$ws.Range("E8").Value = -0.046
$ws.Range("F8").Value = -0.08560000000000001
$ws.Range("K8").Value = 48.1
$ws.Range("L8").Value = 0.19457928802589
$ws.Range("M8").Value = 20.69
$ws.Range("N8").Value = 0.08980034722222222
$ws.Range("O8").Value = 0.4301455301455301
$ws.Range("P8").Value = 15.7
$ws.Range("Q8").Value = 0.0681423611111111
$ws.Range("R8").Value = 0.3264033264033264
$ws.Range("S8").Value = 4.989999999999998
$ws.Range("T8").Value = 0.2411793136781053
$ws.Range("U8").Value = 18.1
$ws.Range("V8").Value = 0.07855902777777778
$ws.Range("W8").Value = 0.06155618121320706
$ws.Range("X8").Value = 0.09397867646572637
$ws.Range("Y8").Value = -0.03242249525251931
$ws.Range("Z8").Value = 0.0934169752853148
$ws.Range("AB8").Value = 0.04091607294513412
$ws.Range("AC8").Value = -0.04091607294513412
$ws.Range("AD8").Value = 2379.4
$ws.Range("AF8").Value = 2379.4
$ws.Range("AG8").Value = 2361.3
$ws.Range("AH8").Value = 0.9117173729787723
$ws.Range("AI8").Value = 0.7471347379658996
$ws.Range("AJ8").Value = 0.9111008218543812
$ws.Range("AK8").Value = 0.7456893829343776

Write-Output "Applied capital structure database update"
